$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Статистика по годам")

function Swap-Rows($ws, $r1, $r2) {
    for ($col = 1; $col -le 5; $col++) {
        $v1 = $ws.Cells.Item($r1, $col).Value2
        $v2 = $ws.Cells.Item($r2, $col).Value2
        $ws.Cells.Item($r1, $col).Value2 = $v2
        $ws.Cells.Item($r2, $col).Value2 = $v1
    }
}

# Swap data rows 3 and 4 (years 2008/2009)
Swap-Rows $ws 3 4

# Swap data rows 6 and 7 (years 2011/2012)
Swap-Rows $ws 6 7
